$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INPUT_SHEET")

# Data rows 2-26: columns A and C hold numeric-looking IDs that must be
# stored as text (matching the existing dataset convention), while B and D
# are naturally non-numeric text ("ONREG-xxxxx" / "ACTIVE").
$rows = @(
    @("100147025","ONREG-23902","64810769121","ACTIVE"),
    @("101487374","ONREG-17938","07782833121","ACTIVE"),
    @("102440213","ONREG-22294","23181610","ACTIVE"),
    @("110062018","ONREG-18131","00300769121","ACTIVE"),
    @("127178114","ONREG-18240","85415732011","ACTIVE"),
    @("139159854","ONREG-25876","09197321011","ACTIVE"),
    @("140389560","ONREG-25888","26262461011","ACTIVE"),
    @("137199660","ONREG-19423","17969657121","ACTIVE"),
    @("100581355","ONREG-19642","45630660011","ACTIVE"),
    @("135142501","ONREG-19769","03757000221","ACTIVE"),
    @("102522400","ONREG-19795","30120160221","ACTIVE"),
    @("109514216","ONREG-19814","12487468121","ACTIVE"),
    @("119002359","ONREG-19836","02550689121","ACTIVE"),
    @("108482050","ONREG-18965","30412792121","ACTIVE"),
    @("163203979","ONREG-23937","15887256121","ACTIVE"),
    @("146457800","ONREG-20446","32958100221","ACTIVE"),
    @("119592056","ONREG-20191","12792592121","ACTIVE"),
    @("164660957","ONREG-12902","12243435121","ACTIVE"),
    @("113526468","ONREG-12961","09212619021","ACTIVE"),
    @("150744017","ONREG-13268","19741512011","ACTIVE"),
    @("131704213","ONREG-13351","33674042011","ACTIVE"),
    @("131478014","ONREG-13415","08509727121","ACTIVE"),
    @("139637466","ONREG-12927","13495060221","ACTIVE"),
    @("142320256","ONREG-13024","07668669121","ACTIVE"),
    @("124700667","ONREG-13036","09297122221","ACTIVE"),
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).NumberFormat = "@"
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,1).Style = "Normal"

    $ws.Cells.Item($r,2).Value = $row[1]

    $ws.Cells.Item($r,3).NumberFormat = "@"
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,3).Style = "Normal"

    $ws.Cells.Item($r,4).Value = $row[3]

    $r = $r + 1
}

# Row 27: appended manually - A27 is a genuine number, B27 repeats the last
# ONREG id as text. No C/D values on this row.
$ws.Cells.Item(27,1).Value = 124700667
$ws.Cells.Item(27,2).Value = "ONREG-13036"

# Restore the view state captured in the saved workbook
$ws.Application.ActiveWindow.ScrollRow = 18
$ws.Range("A25").Select()
